$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "5nt length" results section ---
# Section title
$ws.Range("B24").Value = "5nt length"

# Header row for the new table
$headers = New-Object 'object[,]' 1,4
$headers[0,0] = "Start"
$headers[0,1] = "Stop"
$headers[0,2] = "pval"
$headers[0,3] = "FDR"
$ws.Range("A25:D25").Value = $headers

# Data rows (26-54): Start, Stop, pval, FDR
$data = New-Object 'object[,]' 29,4
  $data[0,0]=209; $data[0,1]=266; $data[0,2]=[double]"7.7811987004626498E-5"; $data[0,3]=[double]"2.33825020948903E-2"
  $data[1,0]=978; $data[1,1]=1001; $data[1,2]=[double]"6.3024163246154898E-3"; $data[1,3]=0.14079089052166199
  $data[2,0]=1028; $data[2,1]=1103; $data[2,2]=[double]"2.2068408411406602E-5"; $data[2,3]=[double]"1.3263113455255399E-2"
  $data[3,0]=1114; $data[3,1]=1132; $data[3,2]=[double]"1.67083740234375E-3"; $data[3,3]=[double]"9.1913743452592303E-2"
  $data[4,0]=1304; $data[4,1]=1329; $data[4,2]=[double]"4.1836539124160098E-3"; $data[4,3]=0.116724861386616
  $data[5,0]=5235; $data[5,1]=5280; $data[5,2]=[double]"1.1953198156854201E-3"; $data[5,3]=[double]"9.1913743452592303E-2"
  $data[6,0]=6403; $data[6,1]=6433; $data[6,2]=[double]"5.7045156426887801E-3"; $data[6,3]=0.14079089052166199
  $data[7,0]=8258; $data[7,1]=8285; $data[7,2]=[double]"6.8774031964892596E-3"; $data[7,3]=0.14761854718178699
  $data[8,0]=8898; $data[8,1]=8919; $data[8,2]=[double]"7.2367882673792404E-3"; $data[8,3]=0.149976198230859
  $data[9,0]=9835; $data[9,1]=9872; $data[9,2]=[double]"6.3250483262643703E-3"; $data[9,3]=0.14079089052166199
  $data[10,0]=9879; $data[10,1]=9902; $data[10,2]=[double]"1.5111633454077101E-3"; $data[10,3]=[double]"9.1913743452592303E-2"
  $data[11,0]=9907; $data[11,1]=9955; $data[11,2]=[double]"2.76078904508382E-3"; $data[11,3]=0.116724861386616
  $data[12,0]=9990; $data[12,1]=10007; $data[12,2]=[double]"1.68228149414063E-3"; $data[12,3]=[double]"9.1913743452592303E-2"
  $data[13,0]=10052; $data[13,1]=10061; $data[13,2]=[double]"3.90625E-3"; $data[13,3]=0.116724861386616
  $data[14,0]=10104; $data[14,1]=10132; $data[14,2]=[double]"1.4833956956863399E-3"; $data[14,3]=[double]"9.1913743452592303E-2"
  $data[15,0]=10182; $data[15,1]=10219; $data[15,2]=[double]"5.1043486463918496E-3"; $data[15,3]=0.13337884941223899
  $data[16,0]=10400; $data[16,1]=10427; $data[16,2]=[double]"9.7234017659598201E-4"; $data[16,3]=[double]"9.1913743452592303E-2"
  $data[17,0]=10499; $data[17,1]=10516; $data[17,2]=[double]"3.9520263671875E-3"; $data[17,3]=0.116724861386616
  $data[18,0]=10557; $data[18,1]=10597; $data[18,2]=[double]"1.47468313388017E-3"; $data[18,3]=[double]"9.1913743452592303E-2"
  $data[19,0]=10684; $data[19,1]=10702; $data[19,2]=[double]"4.0134764726294497E-3"; $data[19,3]=0.116724861386616
  $data[20,0]=11477; $data[20,1]=11491; $data[20,2]=[double]"6.0184740535473898E-3"; $data[20,3]=0.14079089052166199
  $data[21,0]=11558; $data[21,1]=11598; $data[21,2]=[double]"1.9790581624450801E-3"; $data[21,3]=[double]"9.9117829635791299E-2"
  $data[22,0]=14642; $data[22,1]=14663; $data[22,2]=[double]"4.15802001953125E-3"; $data[22,3]=0.116724861386616
  $data[23,0]=15356; $data[23,1]=15385; $data[23,2]=[double]"6.8658081117874299E-4"; $data[23,3]=[double]"9.1913743452592303E-2"
  $data[24,0]=16091; $data[24,1]=16106; $data[24,2]=[double]"1.1884175479869E-3"; $data[24,3]=[double]"9.1913743452592303E-2"
  $data[25,0]=16199; $data[25,1]=16218; $data[25,2]=[double]"4.0134764726294497E-3"; $data[25,3]=0.116724861386616
  $data[26,0]=17032; $data[26,1]=17055; $data[26,2]=[double]"4.2727902670641298E-3"; $data[26,3]=0.116724861386616
  $data[27,0]=17321; $data[27,1]=17348; $data[27,2]=[double]"3.8633549046872698E-3"; $data[27,3]=0.116724861386616
  $data[28,0]=17726; $data[28,1]=17752; $data[28,2]=[double]"2.30634957551956E-3"; $data[28,3]=0.106624314991327
$ws.Range("A26:D54").Value = $data

# Re-apply scientific notation number format to the pval cells that need it
$ws.Range("C26").NumberFormat = "0.00E+00"
$ws.Range("C28").NumberFormat = "0.00E+00"

# Update the visible selection to match the newly added section
$ws.Range("A25:D54").Select() | Out-Null
